$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price/Volume columns hold plain numeric-looking text in the
# source workbook (inline strings), not real numbers. When such text is
# assigned through COM, Excel will silently coerce it into a Number unless the
# destination cell is explicitly formatted as Text first. Collect exactly the
# cells that need a new, numeric-looking value and pre-format only those.

$textRefs = @("D4", "D5", "D6", "D7", "D10", "D11", "D12", "D13", "D14", "D16", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textRefs) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range('D2').Value = '68.342.95'
$ws.Range('E2').Value = '  +4.49%  '

$ws.Range('D3').Value = '3.499.08'
$ws.Range('E3').Value = '  +3.16%  '

$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.29%  '

$ws.Range('D5').Value = '584.48'
$ws.Range('E5').Value = '  +4.45%  '

$ws.Range('D6').Value = '191.18'
$ws.Range('E6').Value = '  +8.85%  '

$ws.Range('D7').Value = '0.633'
$ws.Range('E7').Value = '  +0.13%  '

$ws.Range('D8').Value = '3.474.89'
$ws.Range('E8').Value = '  +2.77%  '

$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('D10').Value = '0.175'
$ws.Range('E10').Value = '  +1.50%  '

$ws.Range('D11').Value = '0.648'
$ws.Range('E11').Value = '  +1.43%  '

$ws.Range('D12').Value = '58.27'
$ws.Range('E12').Value = '  +9.31%  '

$ws.Range('D13').Value = '0.0000280'
$ws.Range('E13').Value = '  +0.56%  '

$ws.Range('D14').Value = '9.52'
$ws.Range('E14').Value = '  +3.36%  '

$ws.Range('D15').Value = '4.022.41'
$ws.Range('E15').Value = '  +2.39%  '

$ws.Range('D16').Value = '19.11'
$ws.Range('E16').Value = '  +4.28%  '

$ws.Range('D17').Value = '3.484.05'
$ws.Range('E17').Value = '  +2.70%  '

$ws.Range('D18').Value = '67.939.08'
$ws.Range('E18').Value = '  +3.96%  '

$ws.Range('D19').Value = '12.24'
$ws.Range('E19').Value = '  +3.34%  '

$ws.Range('E20').Value = '  -0.24%  '

$ws.Range('D21').Value = '1.03'
$ws.Range('E21').Value = '  +2.46%  '

$ws.Range('D22').Value = '488.05'
$ws.Range('E22').Value = '  -0.18%  '

$ws.Range('D23').Value = '5.64'
$ws.Range('E23').Value = '  +14.56%  '

$ws.Range('D24').Value = '17.23'
$ws.Range('E24').Value = '  +21.14%  '

$ws.Range('D25').Value = '4.37'
$ws.Range('E25').Value = '  +6.27%  '

$ws.Range('D26').Value = '89.98'
$ws.Range('E26').Value = '  +0.96%  '

$ws.Range('D27').Value = '3.01'
$ws.Range('E27').Value = '  +3.48%  '

$ws.Range('D28').Value = '11.01'
$ws.Range('E28').Value = '  +2.74%  '

$ws.Range('D29').Value = '9.14'
$ws.Range('E29').Value = '  +4.91%  '

$ws.Range('D30').Value = '31.54'
$ws.Range('E30').Value = '  +0.28%  '

$ws.Range('D31').Value = '7.43'
$ws.Range('E31').Value = '  +13.48%  '

$ws.Range('D32').Value = '610.31'
$ws.Range('E32').Value = '  +6.24%  '

$ws.Range('D33').Value = '11.89'
$ws.Range('E33').Value = '  +3.58%  '

$ws.Range('D34').Value = '64.67'
$ws.Range('E34').Value = '  +3.44%  '

$ws.Range('D35').Value = '0.113'
$ws.Range('E35').Value = '  +4.57%  '

$ws.Range('E36').Value = '  +0.03%  '

$ws.Range('D37').Value = '0.147'
$ws.Range('E37').Value = '  +4.64%  '

$ws.Range('D38').Value = '37.29'
$ws.Range('E38').Value = '  +4.20%  '

$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = '0.391'
$ws.Range('E39').Value = '  +4.71%  '

$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0781'
$ws.Range('E40').Value = '  +5.36%  '

$ws.Range('D41').Value = '3.51'
$ws.Range('E41').Value = '  -2.80%  '

$ws.Range('D42').Value = '3.242.58'
$ws.Range('E42').Value = '  +3.78%  '

$ws.Range('D43').Value = '2.92'
$ws.Range('E43').Value = '  +5.08%  '

$ws.Range('D44').Value = '0.0434'
$ws.Range('E44').Value = '  +4.17%  '

$ws.Range('D45').Value = '2.61'
$ws.Range('E45').Value = '  +7.32%  '

$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '3.25'
$ws.Range('E46').Value = '  +2.62%  '

$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '0.136'
$ws.Range('E47').Value = '  +1.31%  '

$ws.Range('D48').Value = '2.74'
$ws.Range('E48').Value = '  +18.70%  '

$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = '8.74'
$ws.Range('E49').Value = '  +3.82%  '

$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  -0.11%  '

$ws.Range('D51').Value = '140.53'
$ws.Range('E51').Value = '  +0.26%  '
